# Finished Week 13 logging
# Update the "R" (Road) row target-depth tallies on both the OFF and DEF
# sheets with this week's cumulative numbers.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 334
$wsOff.Range("C3").Value = 235
$wsOff.Range("D3").Value = 86
$wsOff.Range("E3").Value = 41
$wsOff.Range("F3").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 364
$wsDef.Range("C3").Value = 248
$wsDef.Range("D3").Value = 88
$wsDef.Range("E3").Value = 45
$wsDef.Range("F3").Value = 4
